$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '35.252.10'
$ws.Range("E2").Value = '  -0.61%  '

# Row 3
$ws.Range("D3").Value = '1.902.22'
$ws.Range("E3").Value = '  +1.49%  '

# Row 4
$ws.Range("E4").Value = '  -0.40%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.12'
$ws.Range("E5").Value = '  +1.54%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.692'
$ws.Range("E6").Value = '  +9.60%  '

# Row 7
$ws.Range("E7").Value = '  -0.37%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.84'
$ws.Range("E8").Value = '  -3.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.351'
$ws.Range("E9").Value = '  +5.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.32'
$ws.Range("E10").Value = '  +11.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0726'
$ws.Range("E11").Value = '  +3.27%  '

# Row 12
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("D13").Value = '2.177.38'
$ws.Range("E13").Value = '  +1.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.32'
$ws.Range("E14").Value = '  +3.65%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.708'
$ws.Range("E15").Value = '  +3.33%  '

# Row 16
$ws.Range("D16").Value = '1.903.50'
$ws.Range("E16").Value = '  +2.09%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.84'
$ws.Range("E17").Value = '  +1.44%  '

# Row 18
$ws.Range("D18").Value = '35.247.42'
$ws.Range("E18").Value = '  -0.71%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.39'
$ws.Range("E19").Value = '  +1.90%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0822'
$ws.Range("E20").Value = '  +2.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.01'
$ws.Range("E21").Value = '  -0.75%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.59'
$ws.Range("E22").Value = '  +2.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("E23").Value = '  +0.86%  '

# Row 24
$ws.Range("E24").Value = '  -0.36%  '

# Row 25
$ws.Range("E25").Value = '  +1.88%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +15.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.69'
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.53'
$ws.Range("E28").Value = '  +3.28%  '

# Row 29
$ws.Range("E29").Value = '  +4.69%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.36'
$ws.Range("E30").Value = '  +3.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.984'
$ws.Range("E32").Value = '  +6.45%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.18'
$ws.Range("E33").Value = '  +2.85%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0571'
$ws.Range("E34").Value = '  +1.24%  '

# Row 35
$ws.Range("E35").Value = '  -0.45%  '

# Row 36
$ws.Range("E36").Value = '  +1.00%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  +0.81%  '

# Row 38
$ws.Range("E38").Value = '  -1.25%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.34'
$ws.Range("E39").Value = '  -1.18%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0682'
$ws.Range("E40").Value = '  +15.66%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$ws.Range("E42").Value = '  +2.98%  '

# Row 43
$ws.Range("E43").Value = '  +6.94%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '90.77'
$ws.Range("E44").Value = '  +1.15%  '

# Row 45
$ws.Range("D45").Value = '1.342.16'
$ws.Range("E45").Value = '  -0.90%  '

# Row 46
$ws.Range("E46").Value = '  +2.76%  '

# Row 47
$ws.Range("B47").Value = 'MultiversX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '47.16'
$ws.Range("E47").Value = '  +4.00%  '

# Row 48
$ws.Range("B48").Value = 'Gas'
$ws.Range("C48").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.79'
$ws.Range("E48").Value = '  +1.94%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.42'
$ws.Range("E49").Value = '  -0.45%  '

# Row 50
$ws.Range("E50").Value = '  +1.95%  '

# Row 51
$ws.Range("E51").Value = '  -2.14%  '
